$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: MIDDLENAME (D2) was re-typed with new sample data ---
$ws.Range("D2").Value = "Samplezczxcz"

# --- New student records appended as rows 3-5 ---
# Row 3
$ws.Range("A3").Value = "2018-11111-BN-0"
$ws.Range("B3").Value = "Samplesfsdf"
$ws.Range("C3").Value = "Sadasdasda"
$ws.Range("D3").Value = "Sampleadad"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Sample with Barangay City Municipality and Country "
$ws.Range("G3").Value = "Male"
$ws.Range("H3").Value = "BSIT"
$ws.Range("I3").Value = "2022, 2023 - 2nd"

# Row 4
$ws.Range("A4").Value = "2018-03300-BN-2"
$ws.Range("B4").Value = "Samplesfsfd"
$ws.Range("C4").Value = "Sample Chazxczxczxnge"
$ws.Range("D4").Value = "Samplejkjkhj"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "Sample with Barangay City Municipality and Country "
$ws.Range("G4").Value = "Male"
$ws.Range("H4").Value = "BSIT"
$ws.Range("I4").Value = "2022, 2023 - 2nd"

# Row 5
$ws.Range("A5").Value = "2018-00550-BN-3"
$ws.Range("B5").Value = "Samplekjkkuku"
$ws.Range("C5").Value = "Sample Change"
$ws.Range("D5").Value = "qeqqweqw"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "Sample with Barangay City Municipality and Country "
$ws.Range("G5").Value = "Male"
$ws.Range("H5").Value = "BSIT"
$ws.Range("I5").Value = "2022, 2023 - 2nd"

# --- Column C (FIRSTNAME) widened to fit the longest new entry ---
$ws.Columns("C").ColumnWidth = 21.42

# --- Selection left on the newly filled-in BATCHYEAR column ---
$ws.Range("I2:I5").Select() | Out-Null
